# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# timestamps on the zh-cn and de-de worksheets to reflect the new handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 15:17:09"
$wsZhCn.Range("H2").Value = "2016-03-23 15:17:32"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 15:17:13"
$wsDeDe.Range("H2").Value = "2016-03-23 15:17:39"
